$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.864.97"
$ws.Range("E2").Value = "  +0.50%  "
$ws.Range("D3").Value = "3.006.11"
$ws.Range("E3").Value = "  +2.14%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "380.49"
$ws.Range("E5").Value = "  +5.89%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "106.08"
$ws.Range("E6").Value = "  +0.85%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.547"
$ws.Range("E7").Value = "  -0.23%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.603"
$ws.Range("E9").Value = "  +1.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.97"
$ws.Range("E10").Value = "  +1.58%  "
$ws.Range("E11").Value = "  -0.30%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0850"
$ws.Range("E12").Value = "  +0.49%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.86"
$ws.Range("E13").Value = "  +0.35%  "
$ws.Range("D14").Value = "3.464.31"
$ws.Range("E14").Value = "  +1.85%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.56"
$ws.Range("E15").Value = "  +1.66%  "
$ws.Range("D16").Value = "2.990.26"
$ws.Range("E16").Value = "  +1.94%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.969"
$ws.Range("E17").Value = "  -0.96%  "
$ws.Range("D18").Value = "51.890.55"
$ws.Range("E18").Value = "  +0.74%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.52"
$ws.Range("E19").Value = "  +6.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.48"
$ws.Range("E20").Value = "  +2.22%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.25"
$ws.Range("E21").Value = "  +0.02%  "
$ws.Range("D22").Value = "0.0₃0965"
$ws.Range("E22").Value = "  +1.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.02"
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "264.87"
$ws.Range("E24").Value = "  +0.47%  "
$ws.Range("E25").Value = "  +4.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.47"
$ws.Range("E26").Value = "  +19.48%  "
$ws.Range("E27").Value = "  -1.99%  "
$ws.Range("B28").Value = "Filecoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.56"
$ws.Range("E28").Value = "  +5.33%  "
$ws.Range("B29").Value = "LEO"
$ws.Range("C29").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.17"
$ws.Range("E29").Value = "  -3.86%  "
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "26.28"
$ws.Range("E30").Value = "  -0.63%  "
$ws.Range("B31").Value = "Dai"
$ws.Range("C31").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.106"
$ws.Range("E32").Value = "  -3.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "10.00"
$ws.Range("E33").Value = "  -0.82%  "
$ws.Range("B34").Value = "InjectiveProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "35.00"
$ws.Range("E34").Value = "  -1.20%  "
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "51.55"
$ws.Range("E35").Value = "  +1.41%  "
$ws.Range("E36").Value = "  -3.71%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0442"
$ws.Range("E37").Value = "  +3.65%  "
$ws.Range("E38").Value = "  +0.13%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.13"
$ws.Range("E39").Value = "  -2.59%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.70"
$ws.Range("E40").Value = "  -5.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.55"
$ws.Range("E41").Value = "  +1.99%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.88"
$ws.Range("E42").Value = "  -0.72%  "
$ws.Range("E43").Value = "  +1.26%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "124.40"
$ws.Range("E44").Value = "  +2.97%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "22.30"
$ws.Range("E45").Value = "  -3.66%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.289"
$ws.Range("E46").Value = "  +21.94%  "
$ws.Range("E47").Value = "  -2.84%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.41"
$ws.Range("E48").Value = "  +5.19%  "
$ws.Range("D49").Value = "2.051.77"
$ws.Range("E49").Value = "  -1.85%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.30"
$ws.Range("E50").Value = "  +1.91%  "
$ws.Range("B51").Value = "BEAM"
$ws.Range("C51").Value = "https://coinranking.com/coin/cYYMfXF4u+beam-beam"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0341"
$ws.Range("E51").Value = "  +6.92%  "
